# Generate Report for Handoff
# Updates status text, "latest" timestamps, and narrows the
# status/datetime column widths on all three sheets.

$wb = $excel.ActiveWorkbook

$newStatus       = "Ready for handoff"
$newOverviewDate = "2016-08-29 15:09:47"
$newZhCnHoDate   = "2016-08-29 15:09:42"
$newDeDeHoDate   = "2016-08-29 15:09:47"
$newColWidth     = 17.2159881591797

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $newOverviewDate
$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = $newZhCnHoDate
$wsZhCn.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = $newDeDeHoDate
$wsDeDe.Columns.Item(3).ColumnWidth = $newColWidth
